$wb = $excel.ActiveWorkbook

# Remove the "Desarquivamentos Pendentes" sheet entirely.
$wsRemove = $wb.Worksheets.Item("Desarquivamentos Pendentes")
$wsRemove.Delete()

# Rename "Paineis DARQ" -> "PAINEIS DARQ" (uppercase).
$wsPaineis = $wb.Worksheets.Item("Paineis DARQ")
$wsPaineis.Name = "PAINEIS DARQ"

# Rename "Recolhimento x Eliminacao" -> "RECOLHIMENTO X ELIMINAÇÃO".
$wsRecolhimento = $wb.Worksheets.Item("Recolhimento x Eliminacao")
$wsRecolhimento.Name = "RECOLHIMENTO X ELIMINAÇÃO"
